$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 stays as text label "test" (unchanged value, only the shared-string
# table shrinks behind the scenes once circle/stripe/none are unused).
$ws.Range("D1").Value = "test"

# Replace the old circle/stripe/none text labels in column D with small
# integer group ids (1-8), one id per distinct Type block.
$typeIds = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 2
    14 = 1
    15 = 3
    16 = 3
    17 = 3
    18 = 3
    19 = 3
    20 = 3
    21 = 3
    22 = 3
    23 = 3
    24 = 3
    25 = 3
    26 = 3
    27 = 3
    28 = 3
    29 = 4
    30 = 4
    31 = 4
    32 = 5
    33 = 5
    34 = 6
    35 = 6
    36 = 7
    37 = 7
    38 = 7
    39 = 7
    40 = 7
    41 = 7
    42 = 7
    43 = 7
    44 = 7
    45 = 7
    46 = 7
    47 = 7
    48 = 7
    49 = 8
}

foreach ($row in $typeIds.Keys) {
    $ws.Cells.Item($row, 4).Value = $typeIds[$row]
}

# The no-longer-referenced "circle"/"stripe"/"none" shared strings should be
# dropped from the workbook now that every D-column cell holds a number.

# Window/view state updates captured in the diff.
$excel.Width = 29040
$excel.Height = 16440
$excel.Left = -120
$excel.Top = -120

$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("G43").Select()
